$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Row 111 currently holds the trailing "※4/8..." note row. A new day's
# figures (2020-05-15) are being appended to the table, so insert a blank
# row above the note -- this shifts the note down to row 112 and keeps its
# original formatting/content intact.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new day's data.
$ws.Cells.Item(111, 1).Value2 = 43966
$ws.Cells.Item(111, 2).Value2 = 216
$ws.Cells.Item(111, 3).Value2 = 37290
$ws.Cells.Item(111, 4).Value2 = 54
$ws.Cells.Item(111, 5).Value2 = 7584

# The data table grew by one row, so the sheet's print area needs to grow
# with it (was $A$1:$E$112, now $A$1:$E$113). Replace the existing
# Print_Area defined name with the updated range.
$wb.Names.Item("相談件数!Print_Area").Delete()
$ws.Names.Add("_xlnm.Print_Area", '=相談件数!$A$1:$E$113')

# Leave the selection on the last cell of new data, matching where the
# editor ended up after appending the new row.
$ws.Range("E109").Select()
